# Update latest output (run 153)
$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("B2").Value = 46043.22916666666
$schedule.Range("C2").Value = 5.5
$schedule.Range("D2").Value = 20.79
$schedule.Range("E2").Value = 716.633034
$schedule.Range("F2").Value = 34.47008340548341

$schedule.Range("A4").Value = 46043.9375
$schedule.Range("C4").Value = 4
$schedule.Range("D4").Value = 15.12
$schedule.Range("E4").Value = 501.94574625
$schedule.Range("F4").Value = 33.19746999007937

$schedule.Range("E5").Value = -9.170791500000012
$schedule.Range("F5").Value = -0.2553826649958232

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("E12").Value = "ON"

$detailed.Range("B41").Value = 79.95
$detailed.Range("B42").Value = 73.37

$detailed.Range("B43").Value = 73.19
$detailed.Range("C43").Value = "historical"

$detailed.Range("B44").Value = 90.24016
$detailed.Range("C44").Value = "historical"

$detailed.Range("B45").Value = 77.80888

$detailed.Range("B46").Value = 77.94
$detailed.Range("E46").Value = "OFF"

$detailed.Range("B49").Value = 64.43329

$detailed.Range("B51").Value = 65.84798000000001
$detailed.Range("B52").Value = 64.09674
$detailed.Range("B53").Value = 63.02046
$detailed.Range("B54").Value = 62.74738
$detailed.Range("B55").Value = 63.60849
$detailed.Range("B56").Value = 66.57365
$detailed.Range("B57").Value = 66.58059
$detailed.Range("B58").Value = 66.56741
$detailed.Range("B59").Value = 67.74234
$detailed.Range("B60").Value = 66.8085
$detailed.Range("B61").Value = 78

$detailed.Range("B66").Value = -5.95454
$detailed.Range("B67").Value = -5.77643
$detailed.Range("B68").Value = -6.35686
$detailed.Range("B69").Value = -7.97915
$detailed.Range("B70").Value = -7.7109
$detailed.Range("B71").Value = -9.193770000000001
$detailed.Range("B72").Value = -7.48385
$detailed.Range("B73").Value = -5.74313
$detailed.Range("B74").Value = -8
$detailed.Range("B75").Value = -8
$detailed.Range("B76").Value = -7.18755
$detailed.Range("B77").Value = -6.00877
$detailed.Range("B78").Value = -5.42834
$detailed.Range("B79").Value = -2.54265

$detailed.Range("B81").Value = -10
$detailed.Range("B82").Value = -11.33055
$detailed.Range("B83").Value = -12.01
$detailed.Range("B84").Value = -12.31954
$detailed.Range("B85").Value = -8.655720000000001
$detailed.Range("B86").Value = -4.9802
$detailed.Range("B87").Value = 0
$detailed.Range("B88").Value = 13.68442
$detailed.Range("B89").Value = 50.55371
$detailed.Range("B90").Value = 29.01365
$detailed.Range("B91").Value = 52.97988

$detailed.Range("B93").Value = 53.96296
$detailed.Range("B94").Value = 30.67231
$detailed.Range("B95").Value = 56.98
$detailed.Range("B96").Value = 56.53663
$detailed.Range("B97").Value = 48.29252
